$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = [double]"0.0005856970953741145"
$ws.Range("E2").Value = [double]"0.0005856970953741145"

# Row 3
$ws.Range("D3").Value = [double]"4.442043758037506E-10"
$ws.Range("E3").Value = [double]"4.442043758037506E-10"

# Row 4
$ws.Range("D4").Value = [double]"2.080953595817627E-20"
$ws.Range("E4").Value = [double]"2.080953595817627E-20"

# Row 5
$ws.Range("D5").Value = [double]"9.134247858550399E-25"
$ws.Range("E5").Value = [double]"9.134247858550399E-25"

# Row 6
$ws.Range("D6").Value = [double]"0.002104692312757382"
$ws.Range("E6").Value = [double]"0.002104692312757382"

# Row 7
$ws.Range("D7").Value = [double]"0.9999999999999982"
$ws.Range("E7").Value = [double]"1.77635683940025E-15"

# Row 8
$ws.Range("D8").Value = [double]"0.9999973039803586"
$ws.Range("E8").Value = [double]"2.696019641423852E-06"

# Row 9
$ws.Range("D9").Value = [double]"0.9999999999992655"
$ws.Range("E9").Value = [double]"7.345235530920036E-13"

# Row 10
$ws.Range("D10").Value = [double]"0.6261004294296078"
$ws.Range("E10").Value = [double]"0.3738995705703922"

# Row 11
$ws.Range("D11").Value = [double]"0.9999988350001748"
$ws.Range("E11").Value = [double]"1.164999825165047E-06"
$ws.Range("F11").Value = [double]"0.04709423333406448"
